$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data keeps every Price/Volume(1h) cell (and the two swapped
# Coin/Link cells) as plain text. For cells whose new value parses as a
# number (e.g. "20.00", "0.0860"), Excel would otherwise auto-convert the
# assignment to a numeric value and silently drop the significant trailing
# zeros / formatting, so force the cell to Text format first.

$ws.Range("D2").Value = "52.322.47"
$ws.Range("E2").Value = "  +1.41%  "
$ws.Range("D3").Value = "2.792.68"
$ws.Range("E3").Value = "  +1.99%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "347.31"
$ws.Range("E5").Value = "  +4.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "115.87"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.551"
$ws.Range("E7").Value = "  +3.97%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.592"
$ws.Range("E9").Value = "  +2.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.49"
$ws.Range("E10").Value = "  +2.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0860"
$ws.Range("E11").Value = "  +4.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.00"
$ws.Range("E12").Value = "  -0.81%  "
$ws.Range("E13").Value = "  +1.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.87"
$ws.Range("E14").Value = "  +3.54%  "
$ws.Range("D15").Value = "3.226.31"
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("D16").Value = "2.854.71"
$ws.Range("E16").Value = "  +4.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.891"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").Value = "52.215.22"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.18"
$ws.Range("E19").Value = "  +6.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.30"
$ws.Range("E20").Value = "  +6.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.34"
$ws.Range("E21").Value = "  -3.54%  "
$ws.Range("D22").Value = "0.0₃0980"
$ws.Range("E22").Value = "  +2.09%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.07"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.68"
$ws.Range("E24").Value = "  -2.03%  "
$ws.Range("E25").Value = "  +3.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.81"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.25"
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("E29").Value = "  +1.00%  "
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.45"
$ws.Range("E31").Value = "  -3.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0456"
$ws.Range("E32").Value = "  +31.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.13"
$ws.Range("E33").Value = "  -0.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.76"
$ws.Range("E34").Value = "  +2.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0832"
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  +0.54%  "
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.60"
$ws.Range("E39").Value = "  -4.01%  "
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.60"
$ws.Range("E41").Value = "  +9.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.35"
$ws.Range("E42").Value = "  -2.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "126.84"
$ws.Range("E43").Value = "  -1.46%  "
$ws.Range("E44").Value = "  +1.98%  "
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.31"
$ws.Range("E46").Value = "  -1.81%  "
$ws.Range("D47").Value = "2.055.17"
$ws.Range("E47").Value = "  -1.98%  "
$ws.Range("E48").Value = "  +3.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.963"
$ws.Range("E49").Value = "  +10.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.62"
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.96"
$ws.Range("E51").Value = "  -0.43%  "
